# Update performance document: fill in the new "v1289" benchmark column (G)
# on the "Sponza" and "ComplexMesh" sheets. The column existed (empty) with
# style already applied; the dependent AVG/VAR/T-TEST/ratio formulas in rows
# 12-16 recompute automatically once the raw samples are present.

$wb = $excel.ActiveWorkbook

# ---- Sponza sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("Sponza")
$ws.Select()

$ws.Range("G1").Value = "v1289"

$sponzaValues = @(10176, 10206, 10194, 10139, 10157, 10184, 10201, 10175, 10125, 10187)
for ($i = 0; $i -lt $sponzaValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $sponzaValues[$i]
}

$ws.Range("G15").Select()

# ---- ComplexMesh sheet ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("ComplexMesh")
$ws2.Select()

$ws2.Range("G1").Value = "v1289"

$complexMeshValues = @(7683, 7657, 7648, 7612, 7716, 7690, 7631, 7615, 7645, 7643)
for ($i = 0; $i -lt $complexMeshValues.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 7).Value = $complexMeshValues[$i]
}

$ws2.Range("G15").Select()

# Keep ComplexMesh as the active/selected tab, matching the source workbook.
$ws2.Select()

$excel.Calculate()
